# Apply the edits described in the commit "Added many more features"
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# 1. Title text (appears as the Heading1 title and again later as a bold
#    paragraph) - both occurrences share the same old/new text, so a single
#    replace-all call takes care of both.
Replace-Text "Play Indian Ruby for Free - Review of Game Structure and Special Symbols" "Play Indian Ruby Online for Free"

# 2. "What we like" bullet list
Replace-Text "Low volatility rating" "Traditional structure with three reels and ten paylines"
Replace-Text "Special Wild and Scatter symbols" "Special symbols and features add excitement to the gameplay"
Replace-Text "Respins function for big wins" "Respins function with a 1000x multiplier for big wins"
Replace-Text "High maximum win potential" "Maximum win of up to €10,000"

# 3. "What we don't like" bullet list (order matters: the first bullet's old
#    text becomes the new text for a different concept than before, and the
#    second bullet's text is replaced by what used to be the first bullet's
#    meaning, with wording tweaked to include "are").
Replace-Text "Highly lucrative wins difficult to achieve" "Average RTP percentage"
Replace-Text "Traditional structure not for everyone" "Highly lucrative wins are difficult to achieve"

# 4. Meta description (italic paragraph)
Replace-Text "Read our review of Indian Ruby slot game, with its traditional structure and special symbols. Play for free and score big with Respins feature." "Read our review of Indian Ruby and play this exciting slot game for free."
